$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the last existing data row (A4) into the new row A5
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# Fill in the new row of time-tracking data
$ws.Range("A5").Value = 42723
$ws.Range("B5").Value = "3h"
$ws.Range("C5").Value = "Mit Code experimentiert, neues Leeres Projekt erstellt, Designskizze angefangen"

# Update selection to match the author's final cursor position
$ws.Range("C6").Select()
